$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team-specific transition-matrix values (NIU_A), updated per the latest
# simulation run. Logic for *using* this data isn't wired up yet -- this
# just refreshes the raw probabilities in place.
$updates = @(
    @{ C = "B2"; V = 0.2513812154696133 },
    @{ C = "C2"; V = 0.4419889502762431 },
    @{ C = "J2"; V = 0.02762430939226519 },
    @{ C = "P2"; V = 0.1685082872928177 },
    @{ C = "S2"; V = 0.1104972375690608 },
    @{ C = "B3"; V = 0.005917159763313609 },
    @{ C = "C3"; V = 0.02958579881656805 },
    @{ C = "J3"; V = 0.01775147928994083 },
    @{ C = "P3"; V = 0.7751479289940828 },
    @{ C = "S3"; V = 0.1715976331360947 },
    @{ C = "J4"; V = 0.1320754716981132 },
    @{ C = "P4"; V = 0.6037735849056604 },
    @{ C = "S4"; V = 0.2641509433962264 },
    @{ C = "B6"; V = 0.0546218487394958 },
    @{ C = "D6"; V = 0.01260504201680672 },
    @{ C = "F6"; V = 0.07142857142857142 },
    @{ C = "J6"; V = 0.2563025210084033 },
    @{ C = "O6"; V = 0.02521008403361345 },
    @{ C = "Q6"; V = 0.1092436974789916 },
    @{ C = "R6"; V = 0.07983193277310924 },
    @{ C = "S6"; V = 0.3907563025210084 },
    @{ C = "B7"; V = 0.1203703703703704 },
    @{ C = "D7"; V = 0.02777777777777778 },
    @{ C = "F7"; V = 0.07870370370370371 },
    @{ C = "J7"; V = 0.1203703703703704 },
    @{ C = "O7"; V = 0.02777777777777778 },
    @{ C = "Q7"; V = 0.125 },
    @{ C = "R7"; V = 0.09259259259259259 },
    @{ C = "S7"; V = 0.4074074074074074 },
    @{ C = "B8"; V = 0.1031746031746032 },
    @{ C = "D8"; V = 0.01984126984126984 },
    @{ C = "E8"; V = 0.001984126984126984 },
    @{ C = "F8"; V = 0.06547619047619048 },
    @{ C = "J8"; V = 0.09126984126984126 },
    @{ C = "O8"; V = 0.01388888888888889 },
    @{ C = "Q8"; V = 0.1646825396825397 },
    @{ C = "R8"; V = 0.125 },
    @{ C = "S8"; V = 0.4146825396825397 },
    @{ C = "B9"; V = 0.0966183574879227 },
    @{ C = "D9"; V = 0.03381642512077294 },
    @{ C = "F9"; V = 0.05797101449275362 },
    @{ C = "J9"; V = 0.1449275362318841 },
    @{ C = "O9"; V = 0.04347826086956522 },
    @{ C = "Q9"; V = 0.1256038647342995 },
    @{ C = "R9"; V = 0.07729468599033816 },
    @{ C = "S9"; V = 0.4202898550724637 },
    @{ C = "B10"; V = 0.1248012718600954 },
    @{ C = "D10"; V = 0.02305246422893482 },
    @{ C = "E10"; V = 0.000794912559618442 },
    @{ C = "F10"; V = 0.07233704292527822 },
    @{ C = "J10"; V = 0.1017488076311606 },
    @{ C = "O10"; V = 0.02066772655007949 },
    @{ C = "Q10"; V = 0.2329093799682035 },
    @{ C = "R10"; V = 0.06677265500794913 },
    @{ C = "S10"; V = 0.3569157392686804 },
    @{ C = "G11"; V = 0.15625 },
    @{ C = "J11"; V = 0.06818181818181818 },
    @{ C = "K11"; V = 0.2159090909090909 },
    @{ C = "L11"; V = 0.5340909090909091 },
    @{ C = "S11"; V = 0.02556818181818182 },
    @{ C = "G12"; V = 0.71875 },
    @{ C = "J12"; V = 0.25 },
    @{ C = "L12"; V = 0.01041666666666667 },
    @{ C = "S12"; V = 0.02083333333333333 },
    @{ C = "F13"; V = 0.02127659574468085 },
    @{ C = "G13"; V = 0.5957446808510638 },
    @{ C = "J13"; V = 0.2978723404255319 },
    @{ C = "S13"; V = 0.0851063829787234 },
    @{ C = "F15"; V = 0.02415458937198068 },
    @{ C = "H15"; V = 0.178743961352657 },
    @{ C = "I15"; V = 0.06763285024154589 },
    @{ C = "J15"; V = 0.3091787439613526 },
    @{ C = "K15"; V = 0.05314009661835749 },
    @{ C = "M15"; V = 0.004830917874396135 },
    @{ C = "O15"; V = 0.05797101449275362 },
    @{ C = "S15"; V = 0.3043478260869565 },
    @{ C = "F16"; V = 0.009302325581395349 },
    @{ C = "H16"; V = 0.2046511627906977 },
    @{ C = "I16"; V = 0.1023255813953488 },
    @{ C = "J16"; V = 0.3767441860465116 },
    @{ C = "K16"; V = 0.1023255813953488 },
    @{ C = "M16"; V = 0.0186046511627907 },
    @{ C = "O16"; V = 0.03255813953488372 },
    @{ C = "S16"; V = 0.1534883720930233 },
    @{ C = "F17"; V = 0.01978021978021978 },
    @{ C = "H17"; V = 0.156043956043956 },
    @{ C = "I17"; V = 0.0989010989010989 },
    @{ C = "J17"; V = 0.4175824175824176 },
    @{ C = "K17"; V = 0.1098901098901099 },
    @{ C = "M17"; V = 0.02417582417582418 },
    @{ C = "N17"; V = 0.002197802197802198 },
    @{ C = "O17"; V = 0.04615384615384616 },
    @{ C = "S17"; V = 0.1252747252747253 },
    @{ C = "F18"; V = 0.01492537313432836 },
    @{ C = "H18"; V = 0.2537313432835821 },
    @{ C = "I18"; V = 0.07960199004975124 },
    @{ C = "J18"; V = 0.3930348258706468 },
    @{ C = "K18"; V = 0.1194029850746269 },
    @{ C = "M18"; V = 0.01492537313432836 },
    @{ C = "O18"; V = 0.05472636815920398 },
    @{ C = "S18"; V = 0.06965174129353234 },
    @{ C = "F19"; V = 0.01174743024963289 },
    @{ C = "H19"; V = 0.2217327459618209 },
    @{ C = "I19"; V = 0.08370044052863436 },
    @{ C = "J19"; V = 0.342143906020558 },
    @{ C = "K19"; V = 0.1204111600587372 },
    @{ C = "M19"; V = 0.02349486049926578 },
    @{ C = "N19"; V = 0.0007342143906020558 },
    @{ C = "O19"; V = 0.05800293685756241 },
    @{ C = "S19"; V = 0.1380323054331865 }

)

foreach ($u in $updates) {
    $ws.Range($u.C).Value = $u.V
}
